$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every value (prices, percentages, dates) as plain text
# (inline strings), not as numbers. To stop Excel from "smart" auto-converting our
# assigned strings into numeric/percentage values (which would lose formatting such
# as trailing zeros, thousands separators and the literal "%" text), we first force
# each target cell to Text number format, then assign the literal string value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.69%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.40%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.150"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.70%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07784"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.99%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.520"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.67%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "14.11%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.566"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.12%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-6.40%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1959"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.88%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04708"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.82%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09351"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.19%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1042"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.67%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001255"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-5.73%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04174"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.85%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005829"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.92%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2,022.23%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.00%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3391"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.38%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.030"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.37%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1346"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.49%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.3039"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.48%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001274"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.42%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003989"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.18%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001353"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.22%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02587"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-3.18%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05939"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.67%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "75.51%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007924"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.90%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1415"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.76%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008420"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "9.59%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008344"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.00%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3119"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.33%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007676"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "9.86%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.11%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.91%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002622"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-34.44%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.11%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.11%"
